$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source inlineStr formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.774.88"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.679.02"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "601.12"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "167.90"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.547"
$ws.Range("D9").Value = "2.678.88"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "27.93"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "3.171.74"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "67.651.93"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "2.685.22"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "11.75"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "7.88"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").Value = "365.81"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -5.14%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "71.02"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("D27").Value = "10.21"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "2.830.28"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "0.0000103"
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "554.27"
$ws.Range("E31").Value = "  -8.24%  "
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  -4.09%  "
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "1.56"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("D38").Value = "19.53"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "155.43"
$ws.Range("E39").Value = "  -4.80%  "
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("D43").Value = "17.95"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  -7.72%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "40.41"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("E47").Value = "  -5.95%  "
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "153.87"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -4.00%  "

# Restore default (unstyled) appearance for the price column, matching original workbook styling
$ws.Range("D2:D51").Style = "Normal"
